$wb = $excel.ActiveWorkbook

# --- Remove the "coden" field from the "fields" table (sheet "fields") ---
$fieldsWs = $wb.Worksheets.Item("fields")
$fieldsWs.Rows(6).Delete()

# --- Remove the corresponding "coden" mapping row from the "mappings" table (sheet "mappings") ---
$mappingsWs = $wb.Worksheets.Item("mappings")
$mappingsWs.Rows(107).Delete()
